# Update the Lgi1-Adam22 LR-pairs sheet with new TPM-derived values.
# The sending-cluster set now includes "ECs" in addition to "MuSCs", so the
# table grows from 5 data rows (MuSCs -> every target cluster) to 10 data
# rows (ECs -> every target cluster, then MuSCs -> every target cluster).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs",   "Lgi1", "Adam22", "ECs",              1, 0.3333333333333333, 0.000409, 0.001227, 0.164521319388576, 0.228024530756365, 3, 1, 1.507281333333333,  4.521844, 0.2777659256678391, 0.2908185693226504, 0.0006164780653333333, 0.005548302588,      0.04569841657206202, 0.06631376780503477),
    @("ECs",   "Lgi1", "Adam22", "FAPs",             1, 0.3333333333333333, 0.000409, 0.001227, 0.164521319388576, 0.228024530756365, 3, 1, 1.225632666666667,  3.676898, 0.2258629392248442, 0.2364765825413957, 0.0005012837606666666, 0.004511553846,      0.03715926876225312, 0.05392246176887056),
    @("ECs",   "Lgi1", "Adam22", "Inflammatory-Mac", 1, 0.3333333333333333, 0.000409, 0.001227, 0.164521319388576, 0.228024530756365, 3, 1, 0.9373149999999999, 2.811945, 0.1727309712259096, 0.1808478624901656, 0.000383361835,         0.003450256515,      0.0284179272853568,  0.04123774898261163),
    @("ECs",   "Lgi1", "Adam22", "MuSCs",            1, 0.3333333333333333, 0.000409, 0.001227, 0.164521319388576, 0.228024530756365, 2, 1, 0.730656,           1.461312, 0.1346472856105345, 0.0939830443451877, 0.000298838304,         0.001793029824,      0.02215234908073557, 0.02143043958586607),
    @("ECs",   "Lgi1", "Adam22", "Resolving-Mac",    1, 0.3333333333333333, 0.000409, 0.001227, 0.164521319388576, 0.228024530756365, 3, 1, 1.025559333333333,  3.076678, 0.1889928782708727, 0.1978739413006007, 0.0004194537673333334, 0.003775083906000001, 0.03109335768816851, 0.04512011261398199),
    @("MuSCs", "Lgi1", "Adam22", "ECs",              1, 0.5,                0.002077, 0.004154, 0.8354786806114239, 0.771975469243635, 3, 1, 1.507281333333333,  4.521844, 0.2777659256678391, 0.2908185693226504, 0.003130623329333333,  0.018783739976,      0.2320675090957771,  0.2245048015176157),
    @("MuSCs", "Lgi1", "Adam22", "FAPs",             1, 0.5,                0.002077, 0.004154, 0.8354786806114239, 0.771975469243635, 3, 1, 1.225632666666667,  3.676898, 0.2258629392248442, 0.2364765825413957, 0.002545639048666666,  0.015273834292,      0.188703670462591,   0.1825541207725251),
    @("MuSCs", "Lgi1", "Adam22", "Inflammatory-Mac", 1, 0.5,                0.002077, 0.004154, 0.8354786806114239, 0.771975469243635, 3, 1, 0.9373149999999999, 2.811945, 0.1727309712259096, 0.1808478624901656, 0.001946803255,         0.01168081953,       0.1443130439405527,  0.139610113507554),
    @("MuSCs", "Lgi1", "Adam22", "MuSCs",            1, 0.5,                0.002077, 0.004154, 0.8354786806114239, 0.771975469243635, 2, 1, 0.730656,           1.461312, 0.1346472856105345, 0.0939830443451877, 0.001517572512,         0.006070290048,      0.1124949365297989,  0.07255260475932164),
    @("MuSCs", "Lgi1", "Adam22", "Resolving-Mac",    1, 0.5,                0.002077, 0.004154, 0.8354786806114239, 0.771975469243635, 3, 1, 1.025559333333333,  3.076678, 0.1889928782708727, 0.1978739413006007, 0.002130086735333334,  0.012780520412,      0.1578995205827042,  0.1527538286866187)
)

$rowIndex = 2
foreach ($row in $data) {
    for ($col = 1; $col -le $row.Length; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $row[$col - 1]
    }
    $rowIndex++
}
